$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 12351
$ws.Range("I40").Value = 5700
$ws.Range("K40").Value = 5700
$ws.Range("M40").Value = -5525

$ws.Range("H62").Value = 4282.6
$ws.Range("I62").Value = 1677.5
$ws.Range("K62").Value = 1677.5
$ws.Range("M62").Value = -1053.5

$ws.Range("H65").Value = 4282.6
$ws.Range("I65").Value = 1677.5
$ws.Range("K65").Value = 8387.5
$ws.Range("M65").Value = -5267.5

$ws.Range("H70").Value = 89333.75
$ws.Range("J70").Value = 7429.4287
$ws.Range("L70").Value = 22288.2861
$ws.Range("N70").Value = -22828.2861

$ws.Range("H73").Value = 89333.75
$ws.Range("J73").Value = 7429.4287
$ws.Range("L73").Value = 22288.2861
$ws.Range("N73").Value = -24160.2861

$ws.Range("H86").Value = 4129.6
$ws.Range("J86").Value = 4210.5
$ws.Range("L86").Value = 4210.5
$ws.Range("N86").Value = -6456.5

$ws.Range("H89").Value = 4129.6
$ws.Range("J89").Value = 4210.5
$ws.Range("L89").Value = 21052.5
$ws.Range("N89").Value = -32284.5

$ws.Range("H96").Value = 166.33333
$ws.Range("I96").Value = 99.5
$ws.Range("J96").Value = 300
$ws.Range("K96").Value = 298.5
$ws.Range("L96").Value = 900
$ws.Range("M96").Value = 1074.5
$ws.Range("N96").Value = -3646

$ws.Range("H106").Value = 5176.6294
$ws.Range("I106").Value = 2344.6155
$ws.Range("K106").Value = 2344.6155
$ws.Range("M106").Value = -1713.6155

$ws.Range("H113").Value = 7167.9
$ws.Range("I113").Value = 9495
$ws.Range("J113").Value = 6909.3335
$ws.Range("K113").Value = 9495
$ws.Range("L113").Value = 6909.3335
$ws.Range("M113").Value = -6241
$ws.Range("N113").Value = -13417.3335

$ws.Range("H132").Value = 1611.8937
$ws.Range("I132").Value = 1307.5682
$ws.Range("K132").Value = 3922.7046
$ws.Range("M132").Value = -1392.7046

$ws.Range("H138").Value = 2986.5
$ws.Range("I138").Value = 1597.8276
$ws.Range("J138").Value = 4478.037
$ws.Range("K138").Value = 4793.4828
$ws.Range("L138").Value = 13434.111
$ws.Range("M138").Value = 346.5172000000002
$ws.Range("N138").Value = -23714.111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9164.583000000001
$ws.Range("I2").Value = 906.8182
$ws.Range("J2").Value = 100000
$ws.Range("K2").Value = 906.8182
$ws.Range("L2").Value = 100000
$ws.Range("M2").Value = -793.8182
$ws.Range("N2").Value = -100226

$ws.Range("H32").Value = 13160116
$ws.Range("I32").Value = 14086363
$ws.Range("K32").Value = 14086363
$ws.Range("M32").Value = -14086076

$ws.Range("H61").Value = 9703.755999999999
$ws.Range("I61").Value = 6540.8057
$ws.Range("K61").Value = 6540.8057
$ws.Range("M61").Value = -6328.8057

$ws.Range("H97").Value = 1067.4375
$ws.Range("I97").Value = 1071.9333
$ws.Range("K97").Value = 1071.9333
$ws.Range("M97").Value = -575.9332999999999

$ws.Range("H102").Value = 5996.2856
$ws.Range("I102").Value = 5995.6665
$ws.Range("K102").Value = 5995.6665
$ws.Range("M102").Value = -4373.6665

$ws.Range("H116").Value = 9164.583000000001
$ws.Range("I116").Value = 906.8182
$ws.Range("J116").Value = 100000
$ws.Range("K116").Value = 906.8182
$ws.Range("L116").Value = 100000
$ws.Range("M116").Value = 1387.1818
$ws.Range("N116").Value = -104588

$ws.Range("H132").Value = 2369.16
$ws.Range("I132").Value = 1853.591
$ws.Range("J132").Value = 6150
$ws.Range("K132").Value = 5560.772999999999
$ws.Range("L132").Value = 18450
$ws.Range("M132").Value = -3030.772999999999
$ws.Range("N132").Value = -23510

$ws.Range("H136").Value = 9703.755999999999
$ws.Range("I136").Value = 6540.8057
$ws.Range("K136").Value = 19622.4171
$ws.Range("M136").Value = -17072.4171

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9164.583000000001
$ws.Range("I3").Value = 906.8182
$ws.Range("J3").Value = 100000
$ws.Range("K3").Value = 906.8182
$ws.Range("L3").Value = 100000
$ws.Range("M3").Value = -792.8182
$ws.Range("N3").Value = -100228

$ws.Range("H86").Value = 4865.643
$ws.Range("J86").Value = 6934.1665
$ws.Range("L86").Value = 6934.1665
$ws.Range("N86").Value = -9180.166499999999

$ws.Range("H89").Value = 4865.643
$ws.Range("J89").Value = 6934.1665
$ws.Range("L89").Value = 34670.8325
$ws.Range("N89").Value = -45902.8325

$ws.Range("H105").Value = 9348.138999999999
$ws.Range("I105").Value = 8444.134
$ws.Range("K105").Value = 8444.134
$ws.Range("M105").Value = -6697.134

$ws.Range("H134").Value = 1932.8667
$ws.Range("I134").Value = 1074.96
$ws.Range("K134").Value = 3224.88
$ws.Range("M134").Value = -689.8800000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 43179.8
$ws.Range("I82").Value = 35000
$ws.Range("J82").Value = 45224.75
$ws.Range("K82").Value = 35000
$ws.Range("L82").Value = 45224.75
$ws.Range("M82").Value = -34639
$ws.Range("N82").Value = -45946.75

$ws.Range("H85").Value = 43179.8
$ws.Range("I85").Value = 35000
$ws.Range("J85").Value = 45224.75
$ws.Range("K85").Value = 35000
$ws.Range("L85").Value = 45224.75
$ws.Range("M85").Value = -33752
$ws.Range("N85").Value = -47720.75

$ws.Range("H134").Value = 6649.7915
$ws.Range("I134").Value = 4582.091
$ws.Range("K134").Value = 13746.273
$ws.Range("M134").Value = -11211.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 7146.1665
$ws.Range("I69").Value = 2388.5
$ws.Range("K69").Value = 7165.5
$ws.Range("M69").Value = -6354.5

$ws.Range("H72").Value = 7146.1665
$ws.Range("I72").Value = 2388.5
$ws.Range("K72").Value = 21496.5
$ws.Range("M72").Value = -17440.5

$ws.Range("H76").Value = 13007.5
$ws.Range("J76").Value = 13007.5
$ws.Range("L76").Value = 39022.5
$ws.Range("N76").Value = -39788.5

$ws.Range("H79").Value = 13007.5
$ws.Range("J79").Value = 13007.5
$ws.Range("L79").Value = 39022.5
$ws.Range("N79").Value = -41674.5

$ws.Range("H131").Value = 6265390.5
$ws.Range("I131").Value = 9616832
$ws.Range("J131").Value = 5118845
$ws.Range("K131").Value = 28850496
$ws.Range("L131").Value = 15356535
$ws.Range("M131").Value = -28845456
$ws.Range("N131").Value = -15366615

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3032857.2
$ws.Range("I102").Value = 3511351
$ws.Range("K102").Value = 3511351
$ws.Range("M102").Value = -3509729

$ws.Range("H132").Value = 308477.62
$ws.Range("I132").Value = 348145.4
$ws.Range("K132").Value = 1044436.2
$ws.Range("M132").Value = -1041906.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9439.916999999999
$ws.Range("I132").Value = 5097.2
$ws.Range("J132").Value = 12541.857
$ws.Range("K132").Value = 15291.6
$ws.Range("L132").Value = 37625.571
$ws.Range("M132").Value = -12761.6
$ws.Range("N132").Value = -42685.571

$ws.Range("H136").Value = 7156.125
$ws.Range("I136").Value = 2849.4
$ws.Range("J136").Value = 14334
$ws.Range("K136").Value = 8548.200000000001
$ws.Range("L136").Value = 43002
$ws.Range("M136").Value = -5998.200000000001
$ws.Range("N136").Value = -48102

$ws.Range("H137").Value = 69999
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 69999
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 69999
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -80199

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 20094.5
$ws.Range("J88").Value = 20094.5
$ws.Range("L88").Value = 20094.5
$ws.Range("N88").Value = -20906.5

$ws.Range("H91").Value = 20094.5
$ws.Range("J91").Value = 20094.5
$ws.Range("L91").Value = 20094.5
$ws.Range("N91").Value = -22902.5

$ws.Range("H122").Value = 10959.667
$ws.Range("I122").Value = 2916.3333
$ws.Range("K122").Value = 8748.999899999999
$ws.Range("M122").Value = -6298.999899999999
